# New crime data collected - update weekly CompStat figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text (volume number + report date range) ---
$ws.Range("A8").Value = "Volume 31   Number  10"
$ws.Range("C9").Value = "Report Covering the Week  3/4/2024  Through  3/10/2024"

# --- Row 15: Rape ---
$ws.Range("L15").Value = -60
$ws.Range("N15").Value = -66.666666666666

# --- Row 16: Robbery ---
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -40
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -16.666666666666
$ws.Range("I16").Value = 27
$ws.Range("J16").Value = 23
$ws.Range("K16").Value = 17.391304347826
$ws.Range("L16").Value = 68.75
$ws.Range("M16").Value = -12.903225806451
$ws.Range("N16").Value = -73.529411764705

# --- Row 17: Fel. Assault ---
$ws.Range("C17").Value = 7
$ws.Range("E17").Value = 133.333333333333
$ws.Range("F17").Value = 26
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = 188.888888888889
$ws.Range("I17").Value = 58
$ws.Range("J17").Value = 36
$ws.Range("K17").Value = 61.111111111111
$ws.Range("L17").Value = 31.818181818181
$ws.Range("M17").Value = 107.142857142857
$ws.Range("N17").Value = -21.621621621621

# --- Row 18: Burglary ---
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 20
$ws.Range("I18").Value = 15
$ws.Range("J18").Value = 15
$ws.Range("L18").Value = -16.666666666666
$ws.Range("M18").Value = -40
$ws.Range("N18").Value = -83.146067415730

# --- Row 19: Gr. Larceny ---
$ws.Range("F19").Value = 16
$ws.Range("G19").Value = 9
$ws.Range("H19").Value = 77.777777777777
$ws.Range("I19").Value = 29
$ws.Range("J19").Value = 32
$ws.Range("K19").Value = -9.375
$ws.Range("L19").Value = -12.121212121212
$ws.Range("M19").Value = 52.631578947368
$ws.Range("N19").Value = -38.297872340425

# --- Row 20: G.L.A. ---
$ws.Range("F20").Value = 6
$ws.Range("H20").Value = -14.285714285714
$ws.Range("I20").Value = 10
$ws.Range("K20").Value = -37.5
$ws.Range("L20").Value = 66.666666666666
$ws.Range("M20").Value = -44.444444444444
$ws.Range("N20").Value = -89.010989010989

# --- Row 21: TOTAL ---
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 11
$ws.Range("E21").Value = 45.454545454545
$ws.Range("F21").Value = 64
$ws.Range("G21").Value = 42
$ws.Range("H21").Value = 52.380952380952
$ws.Range("I21").Value = 142
$ws.Range("J21").Value = 124
$ws.Range("K21").Value = 14.516129032258
$ws.Range("L21").Value = 15.447154471544
$ws.Range("M21").Value = 13.6
$ws.Range("N21").Value = -65.533980582524

# --- Row 22: Transit (C22 flips from a number to the text placeholder "0") ---
$ws.Range("C22").Value = "'0"
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("D22").Value = 2
$ws.Range("E22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E22").Value = -100
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -50
$ws.Range("J22").Value = 5
$ws.Range("K22").Value = -20

# --- Row 23: Housing (C23 flips from text placeholder "0" to a number) ---
$ws.Range("C23").NumberFormat = "#,##0"
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = -28.571428571428
$ws.Range("I23").Value = 14
$ws.Range("J23").Value = 17
$ws.Range("K23").Value = -17.647058823529
$ws.Range("L23").Value = -12.5
$ws.Range("M23").Value = 40

# --- Row 24: Petit Larceny ---
$ws.Range("C24").Value = 8
$ws.Range("D24").Value = 8
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 45
$ws.Range("G24").Value = 31
$ws.Range("H24").Value = 45.161290322580
$ws.Range("I24").Value = 113
$ws.Range("J24").Value = 85
$ws.Range("K24").Value = 32.941176470588
$ws.Range("L24").Value = 26.966292134831
$ws.Range("M24").Value = 91.525423728813

# --- Row 25: Retail Theft (D25,E25 flip from text placeholders to numbers) ---
$ws.Range("C25").Value = 1
$ws.Range("D25").NumberFormat = "#,##0"
$ws.Range("D25").Value = 2
$ws.Range("E25").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E25").Value = -50
$ws.Range("F25").Value = 3
$ws.Range("G25").Value = 6
$ws.Range("H25").Value = -50
$ws.Range("I25").Value = 9
$ws.Range("J25").Value = 11
$ws.Range("K25").Value = -18.181818181818

# --- Row 26: Misd. Assault ---
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = -42.857142857142
$ws.Range("F26").Value = 18
$ws.Range("G26").Value = 22
$ws.Range("H26").Value = -18.181818181818
$ws.Range("I26").Value = 67
$ws.Range("J26").Value = 62
$ws.Range("K26").Value = 8.064516129032
$ws.Range("L26").Value = 6.349206349206
$ws.Range("M26").Value = -4.285714285714

# --- Row 27: UCR Rape* (C27,D27,E27,F27 flip from text placeholders to numbers) ---
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("C27").Value = 1
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("D27").Value = 1
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E27").Value = 0
$ws.Range("F27").NumberFormat = "#,##0"
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -75
$ws.Range("I27").Value = 3
$ws.Range("J27").Value = 6
$ws.Range("K27").Value = -50
$ws.Range("L27").Value = -40

# --- Row 28: Other Sex Crimes (C28,D28,E28 flip from text placeholders to numbers) ---
$ws.Range("C28").NumberFormat = "#,##0"
$ws.Range("C28").Value = 1
$ws.Range("D28").NumberFormat = "#,##0"
$ws.Range("D28").Value = 1
$ws.Range("E28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E28").Value = 0
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = -25
$ws.Range("I28").Value = 8
$ws.Range("J28").Value = 8
$ws.Range("L28").Value = -27.272727272727

# --- Row 29: Shooting Vic. ---
$ws.Range("M29").Value = -66.666666666666
$ws.Range("N29").Value = -90.909090909090

# --- Row 30: Shooting Inc. ---
$ws.Range("M30").Value = -66.666666666666
$ws.Range("N30").Value = -90.909090909090
